$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Fix the product name value on both sheets (added a dash: 341MS -> 341-MS)
$newName = "341-MS-EPP-DB-SAR-REC-NON-RNI-CTPD-SAR-MD-TR-1-ONTIME"
$wsInput.Range("B1").Value = $newName
$wsOutput.Range("B1").Value = $newName

# Update selections on each sheet
$wsInput.Range("B1").Select() | Out-Null
$wsOutput.Range("B1").Select() | Out-Null

# Make ProductLoanOutput the active sheet (activeTab = 1)
$wsOutput.Activate()
